# Updates 1h volume/price snapshot values pulled from coinranking.com
# (GitHub Actions scheduled refresh). Also fixes a rank swap for two
# coin pairs (Litecoin/ShibaInu and Monero/BinanceUSD).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to treat the value as literal text
# (prevents values such as '62.80' or '1.00' from being coerced to
# numbers, matching how this sheet stores every Price/Volume cell as
# a plain string). Resetting the style afterwards clears the resulting
# "quote prefix" formatting flag so no visible style change is left
# behind.
$updates = [ordered]@{
    'D2' = "'" + '25.915.78'
    'E2' = "'" + '  +0.66%  '
    'D3' = "'" + '1.632.41'
    'E3' = "'" + '  +0.27%  '
    'E4' = "'" + '  +0.45%  '
    'D5' = "'" + '214.64'
    'E5' = "'" + '  +0.11%  '
    'E6' = "'" + '  +0.76%  '
    'E7' = "'" + '  +0.42%  '
    'E8' = "'" + '  +0.03%  '
    'D9' = "'" + '0.0631'
    'E9' = "'" + '  -0.03%  '
    'D10' = "'" + '19.59'
    'E10' = "'" + '  +0.47%  '
    'D11' = "'" + '0.0792'
    'E11' = "'" + '  +0.01%  '
    'D12' = "'" + '1.858.58'
    'E12' = "'" + '  +0.32%  '
    'D13' = "'" + '1.632.65'
    'E13' = "'" + '  +0.61%  '
    'E14' = "'" + '  -0.31%  '
    'D15' = "'" + '0.543'
    'E15' = "'" + '  -1.72%  '
    'B16' = "'" + 'Litecoin'
    'C16' = "'" + 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    'D16' = "'" + '62.80'
    'E16' = "'" + '  +0.09%  '
    'B17' = "'" + 'ShibaInu'
    'C17' = "'" + 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D17' = "'" + '0.0₃0754'
    'E17' = "'" + '  -0.77%  '
    'D18' = "'" + '25.907.73'
    'E18' = "'" + '  +0.68%  '
    'D20' = "'" + '192.90'
    'E20' = "'" + '  +0.95%  '
    'E21' = "'" + '  -1.35%  '
    'D22' = "'" + '9.94'
    'E22' = "'" + '  +0.43%  '
    'E23' = "'" + '  -0.23%  '
    'E24' = "'" + '  -0.71%  '
    'B25' = "'" + 'Monero'
    'C25' = "'" + 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D25' = "'" + '143.06'
    'E25' = "'" + '  +0.44%  '
    'B26' = "'" + 'BinanceUSD'
    'C26' = "'" + 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'D26' = "'" + '1.00'
    'E26' = "'" + '  +0.19%  '
    'E27' = "'" + '  +1.88%  '
    'E28' = "'" + '  +0.24%  '
    'E29' = "'" + '  -0.12%  '
    'E30' = "'" + '  +0.08%  '
    'E31' = "'" + '  +0.90%  '
    'D32' = "'" + '3.30'
    'E32' = "'" + '  -0.50%  '
    'E33' = "'" + '  +0.01%  '
    'E34' = "'" + '  -0.18%  '
    'E35' = "'" + '  +2.06%  '
    'E36' = "'" + '  -0.39%  '
    'D37' = "'" + '1.136.58'
    'E37' = "'" + '  -0.10%  '
    'E38' = "'" + '  +1.39%  '
    'E39' = "'" + '  -1.11%  '
    'E40' = "'" + '  +0.57%  '
    'E41' = "'" + '  +0.32%  '
    'D42' = "'" + '0.803'
    'E42' = "'" + '  +0.12%  '
    'D43' = "'" + '99.23'
    'E44' = "'" + '  -1.73%  '
    'D45' = "'" + '1.768.46'
    'E45' = "'" + '  +0.34%  '
    'E46' = "'" + '  +0.83%  '
    'D47' = "'" + '56.26'
    'E47' = "'" + '  +2.05%  '
    'D48' = "'" + '0.0524'
    'E48' = "'" + '  +2.57%  '
    'E49' = "'" + '  +1.58%  '
    'E50' = "'" + '  -0.32%  '
    'E51' = "'" + '  +1.61%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
